$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 435.625
$ws.Range("I5").Value = 768
$ws.Range("K5").Value = 768
$ws.Range("M5").Value = -653

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 520.7222
$ws.Range("I28").Value = 482.75
$ws.Range("K28").Value = 482.75
$ws.Range("M28").Value = 2.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3135.8125
$ws.Range("I40").Value = 2983.3333
$ws.Range("J40").Value = 3227.3
$ws.Range("K40").Value = 2983.3333
$ws.Range("L40").Value = 3227.3
$ws.Range("M40").Value = -2808.3333
$ws.Range("N40").Value = -3577.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 166670690
$ws.Range("I86").Value = 250003570
$ws.Range("K86").Value = 250003570
$ws.Range("M86").Value = -250002447

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 166670690
$ws.Range("I89").Value = 250003570
$ws.Range("K89").Value = 1250017850
$ws.Range("M89").Value = -1250012234

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4415.2856
$ws.Range("I98").Value = 4593.8
$ws.Range("J98").Value = 3969
$ws.Range("K98").Value = 4593.8
$ws.Range("L98").Value = 3969
$ws.Range("M98").Value = -3095.8
$ws.Range("N98").Value = -6965

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2724.1177
$ws.Range("I111").Value = 2697.3635
$ws.Range("J111").Value = 2773.1667
$ws.Range("K111").Value = 8092.0905
$ws.Range("L111").Value = 8319.500100000001
$ws.Range("M111").Value = -5025.0905
$ws.Range("N111").Value = -14453.5001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4415.2856
$ws.Range("I122").Value = 4593.8
$ws.Range("J122").Value = 3969
$ws.Range("K122").Value = 13781.4
$ws.Range("L122").Value = 11907
$ws.Range("M122").Value = -11331.4
$ws.Range("N122").Value = -16807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1106.375
$ws.Range("I2").Value = 1080.1333
$ws.Range("K2").Value = 1080.1333
$ws.Range("M2").Value = -967.1333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 16687500

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1903066.9
$ws.Range("I61").Value = 5718.1665
$ws.Range("K61").Value = 5718.1665
$ws.Range("M61").Value = -5506.1665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 5696.905
$ws.Range("I97").Value = 5696.905
$ws.Range("K97").Value = 5696.905
$ws.Range("M97").Value = -5200.905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 772.5625
$ws.Range("I110").Value = 772.5625
$ws.Range("K110").Value = 772.5625
$ws.Range("M110").Value = 1272.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1106.375
$ws.Range("I116").Value = 1080.1333
$ws.Range("K116").Value = 1080.1333
$ws.Range("M116").Value = 1213.8667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1903066.9
$ws.Range("I136").Value = 5718.1665
$ws.Range("K136").Value = 17154.4995
$ws.Range("M136").Value = -14604.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1106.375
$ws.Range("I3").Value = 1080.1333
$ws.Range("K3").Value = 1080.1333
$ws.Range("M3").Value = -966.1333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1056.3077
$ws.Range("I20").Value = 953.1875
$ws.Range("J20").Value = 1221.3
$ws.Range("K20").Value = 953.1875
$ws.Range("L20").Value = 1221.3
$ws.Range("M20").Value = -706.1875
$ws.Range("N20").Value = -1715.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5896.6665
$ws.Range("I99").Value = 6798
$ws.Range("K99").Value = 6798
$ws.Range("M99").Value = -5300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7750.4814
$ws.Range("I107").Value = 9087.5
$ws.Range("J107").Value = 3930.4285
$ws.Range("K107").Value = 9087.5
$ws.Range("L107").Value = 3930.4285
$ws.Range("M107").Value = -7167.5
$ws.Range("N107").Value = -7770.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2422.5264
$ws.Range("I134").Value = 1901.5385
$ws.Range("J134").Value = 3551.3333
$ws.Range("K134").Value = 5704.6155
$ws.Range("L134").Value = 10653.9999
$ws.Range("M134").Value = -3169.6155
$ws.Range("N134").Value = -15723.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 756318.5
$ws.Range("J141").Value = 738711.2
$ws.Range("L141").Value = 738711.2
$ws.Range("N141").Value = -749071.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 1977.3334
$ws.Range("I40").Value = 3490.3333
$ws.Range("K40").Value = 13961.3332
$ws.Range("M40").Value = -13892.3332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13116.3
$ws.Range("I70").Value = 11589.1875
$ws.Range("K70").Value = 11589.1875
$ws.Range("M70").Value = -11319.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13116.3
$ws.Range("I73").Value = 11589.1875
$ws.Range("K73").Value = 11589.1875
$ws.Range("M73").Value = -10653.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 617.2105
$ws.Range("J97").Value = 939.8
$ws.Range("L97").Value = 939.8
$ws.Range("N97").Value = -1931.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 23988.2
$ws.Range("J136").Value = 23988.2
$ws.Range("L136").Value = 71964.60000000001
$ws.Range("N136").Value = -77064.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1265.9231
$ws.Range("I16").Value = 997.125
$ws.Range("K16").Value = 997.125
$ws.Range("M16").Value = -827.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5219.68
$ws.Range("I46").Value = 25599.75
$ws.Range("K46").Value = 25599.75
$ws.Range("M46").Value = -25411.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3474.6667
$ws.Range("I93").Value = 2274.5
$ws.Range("J93").Value = 5875
$ws.Range("K93").Value = 2274.5
$ws.Range("L93").Value = 5875
$ws.Range("M93").Value = -1026.5
$ws.Range("N93").Value = -8371

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5996
$ws.Range("J136").Value = 5996.25
$ws.Range("L136").Value = 17988.75
$ws.Range("N136").Value = -23088.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 49999.25
$ws.Range("I14").Value = 49999
$ws.Range("K14").Value = 49999
$ws.Range("M14").Value = -49831

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 29907.5
$ws.Range("J55").Value = 29907.5
$ws.Range("L55").Value = 29907.5
$ws.Range("N55").Value = -30461.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2064.6667
$ws.Range("I62").Value = 2597
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 2597
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = -1973
$ws.Range("N62").Value = -2248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 2064.6667
$ws.Range("I65").Value = 2597
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 12985
$ws.Range("L65").Value = 5000
$ws.Range("M65").Value = -9865
$ws.Range("N65").Value = -11240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3091.4
$ws.Range("I126").Value = 3330.5386
$ws.Range("K126").Value = 9991.6158
$ws.Range("M126").Value = -7521.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 57855.277
$ws.Range("I136").Value = 112384.445
$ws.Range("J136").Value = 3326.111
$ws.Range("K136").Value = 337153.335
$ws.Range("L136").Value = 9978.332999999999
$ws.Range("M136").Value = -334603.335
$ws.Range("N136").Value = -15078.333

# Remove cells that no longer exist in the target (structural removals)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N56").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N9").ClearContents()
